# Update automàtic: dades i banners [2026-02-27 21:50]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage cells (column H) need an explicit text format so Excel
# stores the literal "NN%" string instead of auto-converting it to a
# percentage number.
$ws.Range("E2").Value = "2026-02-27 21:48:18"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "55%"
$ws.Range("O2").Value = "5.5 °C"
$ws.Range("E3").Value = "2026-02-27 21:48:21"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "43%"
$ws.Range("O3").Value = "4.3 °C"
$ws.Range("E4").Value = "2026-02-27 21:48:23"
$ws.Range("O4").Value = "9.6 °C"
$ws.Range("E5").Value = "2026-02-27 21:48:26"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "45%"
$ws.Range("N5").Value = "0.5 °C 21:28 TU"
$ws.Range("O5").Value = "4.7 °C"
$ws.Range("E6").Value = "2026-02-27 21:48:28"
$ws.Range("E7").Value = "2026-02-27 21:48:31"
$ws.Range("E8").Value = "2026-02-27 21:48:33"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "67%"
$ws.Range("N8").Value = "8.0 °C 21:11 TU"
$ws.Range("O8").Value = "11.7 °C"
$ws.Range("E9").Value = "2026-02-27 21:48:36"
$ws.Range("E10").Value = "2026-02-27 21:48:38"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "86%"
$ws.Range("N10").Value = "6.4 °C 21:25 TU"
$ws.Range("O10").Value = "10.9 °C"
$ws.Range("E11").Value = "2026-02-27 21:48:40"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "72%"
$ws.Range("O11").Value = "8.5 °C"
$ws.Range("E12").Value = "2026-02-27 21:48:43"
$ws.Range("O12").Value = "10.8 °C"
$ws.Range("E13").Value = "2026-02-27 21:48:45"
$ws.Range("J13").Value = "1025.3 hPa"
$ws.Range("E14").Value = "2026-02-27 21:48:48"
$ws.Range("N14").Value = "5.7 °C 21:21 TU"
$ws.Range("O14").Value = "10.5 °C"
$ws.Range("E15").Value = "2026-02-27 21:48:50"
$ws.Range("E16").Value = "2026-02-27 21:48:52"
$ws.Range("N16").Value = "0.1 °C 21:07 TU"
$ws.Range("O16").Value = "2.6 °C"
$ws.Range("E17").Value = "2026-02-27 21:48:55"
$ws.Range("N17").Value = "4.8 °C 21:13 TU"
$ws.Range("O17").Value = "7.5 °C"
$ws.Range("E18").Value = "2026-02-27 21:48:57"
$ws.Range("E19").Value = "2026-02-27 21:49:00"
$ws.Range("O19").Value = "10.3 °C"
$ws.Range("E20").Value = "2026-02-27 21:49:02"
$ws.Range("E21").Value = "2026-02-27 21:49:05"
$ws.Range("J21").Value = "1024.0 hPa"
$ws.Range("E22").Value = "2026-02-27 21:49:07"
$ws.Range("N22").Value = "-0.4 °C 21:07 TU"
$ws.Range("E23").Value = "2026-02-27 21:49:09"
$ws.Range("O23").Value = "3.6 °C"
$ws.Range("E24").Value = "2026-02-27 21:49:12"
$ws.Range("O24").Value = "10.2 °C"
$ws.Range("E25").Value = "2026-02-27 21:49:14"
$ws.Range("N25").Value = "2.3 °C 21:29 TU"
$ws.Range("O25").Value = "6.0 °C"
$ws.Range("E26").Value = "2026-02-27 21:49:17"
$ws.Range("N26").Value = "6.2 °C 21:29 TU"
$ws.Range("O26").Value = "10.1 °C"
$ws.Range("E27").Value = "2026-02-27 21:49:19"
$ws.Range("N27").Value = "2.5 °C 21:17 TU"
$ws.Range("O27").Value = "5.5 °C"
$ws.Range("E28").Value = "2026-02-27 21:49:21"
$ws.Range("E29").Value = "2026-02-27 21:49:24"
$ws.Range("E30").Value = "2026-02-27 21:49:26"
$ws.Range("E31").Value = "2026-02-27 21:49:29"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "93%"
$ws.Range("E32").Value = "2026-02-27 21:49:31"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "56%"
$ws.Range("E33").Value = "2026-02-27 21:49:34"
$ws.Range("E34").Value = "2026-02-27 21:49:36"
$ws.Range("E35").Value = "2026-02-27 21:49:38"
$ws.Range("J35").Value = "1022.4 hPa"
$ws.Range("E36").Value = "2026-02-27 21:49:41"
$ws.Range("E37").Value = "2026-02-27 21:49:43"
$ws.Range("O37").Value = "8.4 °C"
$ws.Range("E38").Value = "2026-02-27 21:49:46"
$ws.Range("O38").Value = "10.1 °C"
$ws.Range("E39").Value = "2026-02-27 21:49:48"
$ws.Range("N39").Value = "2.3 °C 21:29 TU"
$ws.Range("E40").Value = "2026-02-27 21:49:50"
$ws.Range("J40").Value = "1024.6 hPa"
$ws.Range("O40").Value = "9.0 °C"
$ws.Range("E41").Value = "2026-02-27 21:49:53"
$ws.Range("J41").Value = "1024.5 hPa"
$ws.Range("E42").Value = "2026-02-27 21:49:55"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "90%"
$ws.Range("O42").Value = "11.4 °C"
$ws.Range("E43").Value = "2026-02-27 21:49:57"
$ws.Range("K43").Value = "14.1 MJ/m2"
$ws.Range("E44").Value = "2026-02-27 21:50:00"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "61%"
$ws.Range("E45").Value = "2026-02-27 21:50:02"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "45%"
$ws.Range("E46").Value = "2026-02-27 21:50:04"
